$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row changes ---
# B1: "occupancy schedule" -> "building use type"
$ws.Range("B1").Value = "building use type"

# --- Row 2 changes ---
# B2: was the occupancy-office csv path (string); becomes a numeric factor
$ws.Range("B2").Value = 1.1
# E2: cooling setpoint 28 -> 26
$ws.Range("E2").Value = 26
# F2: emission source "SIA" -> "eu"
$ws.Range("F2").Value = "eu"

# --- Row 3: fill in the previously-empty trailing columns ---
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 23
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = "SIA"

# --- New rows 4 and 5: additional weather file options ---
$ws.Range("A4").Value = "C:\Users\walkerl\Documents\code\sia_380-1-full_version\data\Zweisimmen-hour.epw"
$ws.Range("A5").Value = "C:\Users\walkerl\Documents\code\sia_380-1-full_version\data\Zürich-2070-A1B.epw"

# --- Column width tweaks ---
# Target widths (from the authored file) are 77.08984375 and 14.7265625,
# but this COM host quantizes ColumnWidth to whole-pixel steps, so the
# inputs below are chosen to land on the closest reachable width.
$ws.Columns.Item(1).ColumnWidth = 77
$ws.Columns.Item(5).ColumnWidth = 13.8

# --- Selection moves from A2 to C3 ---
[void]$ws.Range("C3").Select()
